$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The row "grandes regiões e unidades da federação" (row 6, a label-only
# row with no data) is removed. All subsequent rows shift up by one.
$ws.Rows("6:6").Delete()
